# Week 13 logging update
# Appends newly-logged play values to the running per-game log strings on
# the YDS and ST sheets, and updates the week's aggregate numeric totals on
# the OFF, DEF, ST, TURNS and PEN sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# YDS sheet — append this week's individual play yardage logs
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("YDS")

$ws.Range("B2").Value = $ws.Range("B2").Value2 + " 4 2 4 2 4 0 5 7 0 7 5 5 1 7 16 2 2 8 3"
$ws.Range("C2").Value = $ws.Range("C2").Value2 + " 2 10 8 6 2 1 10 1 -1 6 0 7 6 0 3 8 2 0 7 0 -1 15 11 -1 4 9 3 16 4 1 11 25 4 5 3 6 -3 7 7 3 6 10 0 12 3"
$ws.Range("B3").Value = $ws.Range("B3").Value2 + " 15 5 12 11 2 10 4 14 2 5"
$ws.Range("C3").Value = $ws.Range("C3").Value2 + " 12 4 8 9 14 6 13 2 8 17 24 14 4 17 0 6"

# ---------------------------------------------------------------------
# OFF sheet — week 13 offensive totals
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("OFF")

$ws.Range("B2").Value = 6
$ws.Range("C2").Value = 115
$ws.Range("E2").Value = 10
$ws.Range("F2").Value = 61
$ws.Range("G2").Value = 27
$ws.Range("H2").Value = 3
$ws.Range("I2").Value = 6
$ws.Range("J2").Value = 22
$ws.Range("L2").Value = 184
$ws.Range("M2").Value = 120
$ws.Range("O2").Value = 17
$ws.Range("Q2").Value = 347

$ws.Range("B3").Value = 6
$ws.Range("C3").Value = 104
$ws.Range("E3").Value = 20
$ws.Range("F3").Value = 62
$ws.Range("G3").Value = 21
$ws.Range("H3").Value = 29
$ws.Range("I3").Value = 44
$ws.Range("J3").Value = 46
$ws.Range("N3").Value = 22

# ---------------------------------------------------------------------
# DEF sheet — week 13 defensive totals
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("DEF")

$ws.Range("B2").Value = 4
$ws.Range("C2").Value = 169
$ws.Range("D2").Value = 9
$ws.Range("E2").Value = 8
$ws.Range("F2").Value = 40
$ws.Range("G2").Value = 49
$ws.Range("I2").Value = 3
$ws.Range("J2").Value = 23
$ws.Range("L2").Value = 193
$ws.Range("M2").Value = 125
$ws.Range("O2").Value = 14
$ws.Range("P2").Value = 9
$ws.Range("Q2").Value = 386

$ws.Range("B3").Value = 8
$ws.Range("C3").Value = 110
$ws.Range("D3").Value = 7
$ws.Range("E3").Value = 27
$ws.Range("F3").Value = 79
$ws.Range("H3").Value = 20
$ws.Range("I3").Value = 33
$ws.Range("J3").Value = 37
$ws.Range("N3").Value = 10

# ---------------------------------------------------------------------
# ST sheet — week 13 special-teams totals and logs
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ST")

$ws.Range("B2").Value = 39
$ws.Range("D2").Value = 65
$ws.Range("B3").Value = 22

$ws.Range("D3").Value = $ws.Range("D3").Value2 + " 50 47 53 60 41 51"
$ws.Range("D4").Value = $ws.Range("D4").Value2 + " 7 -1 16 0 0 15"
$ws.Range("D5").Value = $ws.Range("D5").Value2 + " 0 0 4"
$ws.Range("B6").Value = $ws.Range("B6").Value2 + " 26 22"

# ---------------------------------------------------------------------
# TURNS sheet — week 13 turnover totals
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("TURNS")

$ws.Range("B2").Value = 4
$ws.Range("D2").Value = 7
$ws.Range("E2").Value = 5
$ws.Range("D3").Value = 6
$ws.Range("E3").Value = 4

# ---------------------------------------------------------------------
# PEN sheet — week 13 penalty totals
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("PEN")

$ws.Range("B2").Value = 11
$ws.Range("B3").Value = 20
$ws.Range("B5").Value = 1
